# Auto-generated edit script: updates FFXIV crafting profit data cells
# across sheets ALC, ARM, BSM, CUL, GSM, LTW, WVR per scheduled-runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 885.9091
$ws.Range("I62").Value = 820.1667
$ws.Range("J62").Value = 964.8
$ws.Range("K62").Value = 820.1667
$ws.Range("L62").Value = 964.8
$ws.Range("M62").Value = -196.1667
$ws.Range("N62").Value = -2212.8
$ws.Range("H65").Value = 885.9091
$ws.Range("I65").Value = 820.1667
$ws.Range("J65").Value = 964.8
$ws.Range("K65").Value = 4100.8335
$ws.Range("L65").Value = 4824
$ws.Range("M65").Value = -980.8334999999997
$ws.Range("N65").Value = -11064
$ws.Range("H108").Value = 70000
$ws.Range("J108").Value = 70000
$ws.Range("L108").Value = 70000
$ws.Range("N108").Value = -77680
$ws.Range("H111").Value = 1069.6786
$ws.Range("I111").Value = 1206.125
$ws.Range("J111").Value = 887.75
$ws.Range("K111").Value = 3618.375
$ws.Range("L111").Value = 2663.25
$ws.Range("M111").Value = -551.375
$ws.Range("N111").Value = -8797.25
$ws.Range("H116").Value = 5153638.5
$ws.Range("I116").Value = 6298291.5
$ws.Range("J116").Value = 2698.9
$ws.Range("K116").Value = 6298291.5
$ws.Range("L116").Value = 2698.9
$ws.Range("M116").Value = -6294849.5
$ws.Range("N116").Value = -9582.9
$ws.Range("H132").Value = 4752.5
$ws.Range("I132").Value = 4991.2583
$ws.Range("J132").Value = 4317.1177
$ws.Range("K132").Value = 14973.7749
$ws.Range("L132").Value = 12951.3531
$ws.Range("M132").Value = -12443.7749
$ws.Range("N132").Value = -18011.3531

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12654.094
$ws.Range("I32").Value = 3446.0344
$ws.Range("J32").Value = 32434.371
$ws.Range("K32").Value = 3446.0344
$ws.Range("L32").Value = 32434.371
$ws.Range("M32").Value = -3159.0344
$ws.Range("N32").Value = -33008.371
$ws.Range("H61").Value = 1096.6216
$ws.Range("I61").Value = 948.1429000000001
$ws.Range("J61").Value = 1291.5
$ws.Range("K61").Value = 948.1429000000001
$ws.Range("L61").Value = 1291.5
$ws.Range("M61").Value = -736.1429000000001
$ws.Range("N61").Value = -1715.5
$ws.Range("H122").Value = 2668.6667
$ws.Range("I122").Value = 1848.8889
$ws.Range("K122").Value = 5546.6667
$ws.Range("M122").Value = -3096.6667
$ws.Range("H123").Value = 50250
$ws.Range("J123").Value = 50250
$ws.Range("L123").Value = 50250
$ws.Range("N123").Value = -60050
$ws.Range("H131").Value = 41357.5
$ws.Range("J131").Value = 41357.5
$ws.Range("L131").Value = 41357.5
$ws.Range("N131").Value = -51437.5
$ws.Range("H136").Value = 1096.6216
$ws.Range("I136").Value = 948.1429000000001
$ws.Range("J136").Value = 1291.5
$ws.Range("K136").Value = 2844.4287
$ws.Range("L136").Value = 3874.5
$ws.Range("M136").Value = -294.4287000000004
$ws.Range("N136").Value = -8974.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 66667828
$ws.Range("I86").Value = 66667828
$ws.Range("K86").Value = 66667828
$ws.Range("M86").Value = -66666705
$ws.Range("H89").Value = 66667828
$ws.Range("I89").Value = 66667828
$ws.Range("K89").Value = 333339140
$ws.Range("M89").Value = -333333524
$ws.Range("H105").Value = 11200
$ws.Range("I105").Value = 1500
$ws.Range("J105").Value = 50000
$ws.Range("K105").Value = 1500
$ws.Range("L105").Value = 50000
$ws.Range("M105").Value = 247
$ws.Range("N105").Value = -53494
$ws.Range("H134").Value = 1386.4333
$ws.Range("I134").Value = 1210.1034
$ws.Range("J134").Value = 6500
$ws.Range("K134").Value = 3630.3102
$ws.Range("L134").Value = 19500
$ws.Range("M134").Value = -1095.3102
$ws.Range("N134").Value = -24570

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 996.5952
$ws.Range("I131").Value = 419.2
$ws.Range("J131").Value = 1177.0312
$ws.Range("K131").Value = 1257.6
$ws.Range("L131").Value = 3531.0936
$ws.Range("M131").Value = 3782.4
$ws.Range("N131").Value = -13611.0936
$ws.Range("H137").Value = 4171.3335
$ws.Range("I137").Value = 676.6667
$ws.Range("J137").Value = 7666
$ws.Range("K137").Value = 2030.0001
$ws.Range("L137").Value = 22998
$ws.Range("M137").Value = 3069.9999
$ws.Range("N137").Value = -33198

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 21626012
$ws.Range("I70").Value = 66670830
$ws.Range("J70").Value = 4496.72
$ws.Range("K70").Value = 66670830
$ws.Range("L70").Value = 4496.72
$ws.Range("M70").Value = -66670560
$ws.Range("N70").Value = -5036.72
$ws.Range("H73").Value = 21626012
$ws.Range("I73").Value = 66670830
$ws.Range("J73").Value = 4496.72
$ws.Range("K73").Value = 66670830
$ws.Range("L73").Value = 4496.72
$ws.Range("M73").Value = -66669894
$ws.Range("N73").Value = -6368.72
$ws.Range("H102").Value = 2274.9092
$ws.Range("I102").Value = 2274.9092
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2274.9092
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -652.9092000000001
$ws.Range("N102").ClearContents()
$ws.Range("H122").Value = 5317.3335
$ws.Range("I122").Value = 4672.0435
$ws.Range("J122").Value = 6801.5
$ws.Range("K122").Value = 14016.1305
$ws.Range("L122").Value = 20404.5
$ws.Range("M122").Value = -11566.1305
$ws.Range("N122").Value = -25304.5
$ws.Range("H126").Value = 2983.3333
$ws.Range("I126").Value = 1633.3334
$ws.Range("K126").Value = 4900.0002
$ws.Range("M126").Value = -2430.0002
$ws.Range("H132").Value = 972.25
$ws.Range("I132").Value = 837.3103599999999
$ws.Range("J132").Value = 2276.6667
$ws.Range("K132").Value = 2511.93108
$ws.Range("L132").Value = 6830.000100000001
$ws.Range("M132").Value = 18.06892000000016
$ws.Range("N132").Value = -11890.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 656.2857
$ws.Range("I16").Value = 656.2857
$ws.Range("K16").Value = 656.2857
$ws.Range("M16").Value = -486.2857
$ws.Range("H61").Value = 1646.4117
$ws.Range("I61").Value = 1769
$ws.Range("J61").Value = 1352.2
$ws.Range("K61").Value = 1769
$ws.Range("L61").Value = 1352.2
$ws.Range("M61").Value = -1567
$ws.Range("N61").Value = -1756.2
$ws.Range("H93").Value = 29412538
$ws.Range("I93").Value = 45455280
$ws.Range("J93").Value = 849.6667
$ws.Range("K93").Value = 45455280
$ws.Range("L93").Value = 849.6667
$ws.Range("M93").Value = -45454032
$ws.Range("N93").Value = -3345.6667
$ws.Range("H113").Value = 1646.4117
$ws.Range("I113").Value = 1769
$ws.Range("J113").Value = 1352.2
$ws.Range("K113").Value = 1769
$ws.Range("L113").Value = 1352.2
$ws.Range("M113").Value = 401
$ws.Range("N113").Value = -5692.2
$ws.Range("H122").Value = 5454.5
$ws.Range("I122").Value = 5482.091
$ws.Range("J122").Value = 5151
$ws.Range("K122").Value = 16446.273
$ws.Range("L122").Value = 15453
$ws.Range("M122").Value = -13996.273
$ws.Range("N122").Value = -20353
$ws.Range("H132").Value = 14541817
$ws.Range("I132").Value = 31262166
$ws.Range("K132").Value = 93786498
$ws.Range("M132").Value = -93783968
$ws.Range("H136").Value = 3692.1836
$ws.Range("I136").Value = 5133.6553
$ws.Range("J136").Value = 1602.05
$ws.Range("K136").Value = 15400.9659
$ws.Range("L136").Value = 4806.15
$ws.Range("M136").Value = -12850.9659
$ws.Range("N136").Value = -9906.15

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 9015
$ws.Range("I100").Value = 19996
$ws.Range("J100").Value = 1171.4286
$ws.Range("K100").Value = 39992
$ws.Range("L100").Value = 2342.8572
$ws.Range("M100").Value = -39451
$ws.Range("N100").Value = -3424.8572
$ws.Range("H107").Value = 1700
$ws.Range("I107").Value = 600
$ws.Range("J107").Value = 5000
$ws.Range("K107").Value = 1800
$ws.Range("L107").Value = 15000
$ws.Range("M107").Value = 120
$ws.Range("N107").Value = -18840
$ws.Range("H123").Value = 43344.316
$ws.Range("J123").Value = 43344.316
$ws.Range("L123").Value = 43344.316
$ws.Range("N123").Value = -53144.316
$ws.Range("H132").Value = 1310.6724
$ws.Range("I132").Value = 648.25
$ws.Range("J132").Value = 7051.6665
$ws.Range("K132").Value = 1944.75
$ws.Range("L132").Value = 21154.9995
$ws.Range("M132").Value = 585.25
$ws.Range("N132").Value = -26214.9995
